$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin -> Bitcoin)
$ws.Range("D2").Value = "26.452.45"
$ws.Range("E2").Value = "  +1.56%  "

# Row 3 (Ethereum -> Ethereum)
$ws.Range("D3").Value = "1.677.08"
$ws.Range("E3").Value = "  +2.39%  "

# Row 4 (TetherUSD -> TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 (BNB -> BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.99"
$ws.Range("E5").Value = "  +2.39%  "

# Row 6 (XRP -> XRP)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5336"
$ws.Range("E6").Value = "  +1.89%  "

# Row 7 (USDC -> USDC)
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 (Cardano -> Cardano)
$ws.Range("E8").Value = "  +4.29%  "

# Row 9 (Dogecoin -> Dogecoin)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06408"
$ws.Range("E9").Value = "  +1.88%  "

# Row 10 (Solana -> Solana)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  +6.39%  "

# Row 11 (TRON -> TRON)
$ws.Range("E11").Value = "  +1.69%  "

# Row 12 (Polkadot -> Polkadot)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.514"
$ws.Range("E12").Value = "  +2.74%  "

# Row 13 (WrappedEther -> WrappedEther)
$ws.Range("D13").Value = "1.674.42"
$ws.Range("E13").Value = "  +2.23%  "

# Row 14 (Polygon -> Polygon)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5592"
$ws.Range("E14").Value = "  +1.48%  "

# Row 15 (ShibaInu -> ShibaInu)
$ws.Range("D15").Value = "0.0₅8325"
$ws.Range("E15").Value = "  +1.85%  "

# Row 16 (Litecoin -> Litecoin)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.69"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17 (WrappedBTC -> WrappedBTC)
$ws.Range("D17").Value = "26.508.36"
$ws.Range("E17").Value = "  +1.79%  "

# Row 18 (Dai -> Dai)
$ws.Range("E18").Value = "  -0.06%  "

# Row 19 (Uniswap -> Uniswap)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.798"
$ws.Range("E19").Value = "  +2.63%  "

# Row 20 (BitcoinCash -> BitcoinCash)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.38"
$ws.Range("E20").Value = "  +2.74%  "

# Row 21 (Avalanche -> Avalanche)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("E21").Value = "  +1.36%  "

# Row 22 (Chainlink -> Chainlink)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.325"
$ws.Range("E22").Value = "  +2.86%  "

# Row 23 (BinanceUSD -> BinanceUSD)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 (Monero -> Stellar)
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1280"
$ws.Range("E24").Value = "  +6.26%  "

# Row 25 (Stellar -> Monero)
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.69"
$ws.Range("E25").Value = "  -3.14%  "

# Row 26 (Cosmos -> Cosmos)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.411"
$ws.Range("E26").Value = "  +0.30%  "

# Row 27 (EthereumClassic -> EthereumClassic)
$ws.Range("E27").Value = "  +3.25%  "

# Row 28 (Toncoin -> Toncoin)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.445"
$ws.Range("E28").Value = "  +4.78%  "

# Row 29 (Hedera -> Hedera)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06282"
$ws.Range("E29").Value = "  +5.44%  "

# Row 30 (PancakeSwap -> PancakeSwap)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.283"
$ws.Range("E30").Value = "  +2.28%  "

# Row 31 (InternetComputer(DFINITY) -> InternetComputer(DFINITY))
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.607"
$ws.Range("E31").Value = "  +5.44%  "

# Row 32 (Filecoin -> Filecoin)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.459"
$ws.Range("E32").Value = "  +2.04%  "

# Row 33 (LidoDAOToken -> LidoDAOToken)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.699"
$ws.Range("E33").Value = "  +3.44%  "

# Row 34 (ARBITRUM -> ARBITRUM)
$ws.Range("E34").Value = "  +3.17%  "

# Row 35 (ImmutableX -> ImmutableX)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6152"
$ws.Range("E35").Value = "  +9.36%  "

# Row 36 (HuobiToken -> HuobiToken)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.423"
$ws.Range("E36").Value = "  +1.14%  "

# Row 37 (MXToken -> MXToken)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("E37").Value = "  +0.88%  "

# Row 38 (VeChain -> FraxShare)
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.145"
$ws.Range("E38").Value = "  +8.10%  "

# Row 39 (FraxShare -> VeChain)
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01631"
$ws.Range("E39").Value = "  +1.16%  "

# Row 40 (Maker -> Maker)
$ws.Range("D40").Value = "1.097.14"
$ws.Range("E40").Value = "  +6.99%  "

# Row 41 (TrustWalletToken -> TrustWalletToken)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8656"
$ws.Range("E41").Value = "  +2.13%  "

# Row 42 (PaxDollar -> PaxDollar)
$ws.Range("E42").Value = "  -0.04%  "

# Row 43 (Quant -> Quant)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +0.55%  "

# Row 44 (RocketPoolETH -> RocketPoolETH)
$ws.Range("D44").Value = "1.822.63"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45 (BabyDogeCoin -> BabyDogeCoin)
$ws.Range("E45").Value = "  +5.06%  "

# Row 46 (Aave -> Aave)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.46"
$ws.Range("E46").Value = "  +4.95%  "

# Row 47 (EnergySwap -> EnergySwap)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.195"
$ws.Range("E47").Value = "  +1.65%  "

# Row 48 (Frax -> Frax)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9997"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49 (Cronos -> Cronos)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05210"
$ws.Range("E49").Value = "  +1.36%  "

# Row 50 (RenderToken -> RenderToken)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.474"
$ws.Range("E50").Value = "  +7.14%  "

# Row 51 (Aptos -> Aptos)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.027"
$ws.Range("E51").Value = "  +2.27%  "
